# Generate Report for Handoff
# Updates the "b.md" rows (row 3) on the Overview, zh-cn and de-de sheets to
# reflect that the file is now "Ready for handoff" (instead of the previous
# "Handed back" state), with refreshed handback file names / timestamps and
# a new error detail message. Also widens column P (Error Detail) on the
# locale sheets so the new, longer message is readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": row 3 corresponds to b.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 18:40:19"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": row 3 corresponds to b.md
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to fit the new message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

$wsZhCn.Range("C3").Value = "Ready for handoff"

# NOTE: assigning the bare string "False" gets auto-interpreted as a Boolean
# by Excel (matching native Excel behaviour). The source file stores this as
# plain text, so force text entry via a leading apostrophe and then clear the
# resulting "quote prefix" style back to Normal.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"

$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-23 18:39:59"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf50354f3436e6be4b24c9d4e8cd0a073b39cb18/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692e452ac70c435f943eee08d353c978929193b1/e2e/b.md."

# ---------------------------------------------------------------------------
# Sheet "de-de": row 3 corresponds to b.md
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P) to fit the new message.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667

$wsDeDe.Range("C3").Value = "Ready for handoff"

# See note above regarding Boolean auto-detection for the literal text "False".
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"

$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-23 18:40:19"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf50354f3436e6be4b24c9d4e8cd0a073b39cb18/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692e452ac70c435f943eee08d353c978929193b1/e2e/b.md."
